$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 19.471208
$ws.Range("D3").Value = 77.54560499999999

$wb.Save()
